$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A46 needs to hold the literal text "2020-07-15" as a shared string,
# not be auto-converted into a date serial number by Excel's input
# parsing. Enter it as a text formula, then convert the formula to a
# plain value in place (copy / paste-values), which yields a normal
# shared-string text cell without adding any new cell style.
$ws.Range("A46").Formula = '="2020-07-15"'
$ws.Range("A46").Copy()
$ws.Range("A46").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B46").Value = 317635
$ws.Range("C46").Value = 369411
$ws.Range("D46").Value = 81411
$ws.Range("E46").Value = 36906
$ws.Range("F46").Value = 29.05
